# Apply the "Anonymize fedcore" update:
#  - rename the "fedcore" column header to "approach" on both sheets
#  - redraw the thin-border box around the merged header cells (B1:D1 / E1:G1)
#    so the box is built from the per-cell edges (top/bottom on the middle
#    cell of the merge, top/bottom/right on the right-most cell) instead of
#    re-using the single all-around border on every cell of the merge
#  - drop the stray empty inline-string placeholder in G5 on the
#    computational_comparison sheet
#
# NOTE on the border construction: building each border combination directly,
# cell-by-cell, via the Borders collection creates a brand-new style record
# (cellXfs entry) for every *intermediate* edge combination a cell passes
# through. When two different cells each transition through the same
# throwaway intermediate combination, that leftover record never gets
# reused/cleaned up, so the workbook ends up with extra unused styles. To
# avoid that we build each final border combination exactly once on a single
# scratch cell (so only the same cell ever mutates through the intermediate
# states) and then stamp the finished look onto every real target cell with
# a formats-only paste.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("quality_comparison")
$ws2 = $wb.Worksheets.Item("computational_comparison")

# --- build the two border looks once, on a scratch cell far outside the
#     used ranges of either sheet ------------------------------------------
$scratch = $ws1.Range("Z100")
$scratch.Clear()

# Look #1: top + bottom only (middle cell of a merged header box)
$scratch.Borders.Item(8).LineStyle = 1    # xlEdgeTop
$scratch.Borders.Item(9).LineStyle = 1    # xlEdgeBottom

$scratch.Copy()
$ws1.Range("C1").PasteSpecial(-4122)      # xlPasteFormats
$ws2.Range("C1").PasteSpecial(-4122)
$ws2.Range("F1").PasteSpecial(-4122)

# Look #2: top + bottom + right (right-most cell of a merged header box) —
# keep mutating the same scratch cell so the change happens in place.
$scratch.Borders.Item(10).LineStyle = 1   # xlEdgeRight

$scratch.Copy()
$ws1.Range("D1").PasteSpecial(-4122)
$ws2.Range("D1").PasteSpecial(-4122)
$ws2.Range("G1").PasteSpecial(-4122)

$scratch.Clear()
$excel.CutCopyMode = 0

# --- anonymize the "fedcore" header text ----------------------------------
$ws1.Range("C2").Value = "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# --- drop the leftover empty placeholder cell -----------------------------
$ws2.Range("G5").ClearContents()
